# Add files via upload
#
# Populate the new "entity / attribute" mini tables that sit below the
# existing "Tabela1" library table (rows 14-17, columns A:C and E:G)
# and turn each new block into its own Excel Table (ListObject), mirroring
# the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values --------------------------------------------------------
# Write order chosen so shared strings are interned in the same order as
# in the source workbook: Coluna1, ENTIDADE, ATRIBUTOS, Relacionamento,
# USUÁRIOS, DADOS, OBRAS, ILIMITADO, INFORMAÇÕES, POSSE.

$ws.Range("C14").Value = "Coluna1"
$ws.Range("A14").Value = "ENTIDADE"
$ws.Range("B14").Value = "ATRIBUTOS"
$ws.Range("E14").Value = "ENTIDADE"
$ws.Range("F14").Value = "ATRIBUTOS"
$ws.Range("G14").Value = "Coluna1"

$ws.Range("G15").Value = "Relacionamento"
$ws.Range("A15").Value = "USUÁRIOS"
$ws.Range("B15").Value = "DADOS"
$ws.Range("E15").Value = "USUÁRIOS"
$ws.Range("F15").Value = "DADOS"

$ws.Range("A16").Value = "OBRAS"
$ws.Range("B16").Value = "ILIMITADO"
$ws.Range("E16").Value = "OBRAS"
$ws.Range("F16").Value = "ILIMITADO"
$ws.Range("G16").Value = "INFORMAÇÕES"

$ws.Range("G17").Value = "POSSE"

# --- Turn the four new ranges into Excel Tables -------------------------
# NOTE: renaming a freshly added ListObject can shuffle the internal
# (alphabetically sorted) ListObjects collection, which would silently
# invalidate any previously captured object reference. So every lookup
# below is re-resolved fresh by range address instead of reusing a
# stashed variable.

function New-NamedTable($ws, $rangeAddr, $name, $style) {
    [void]$ws.ListObjects.Add(1, $ws.Range($rangeAddr), $null, 1)
    $justAdded = $ws.ListObjects.Item($ws.ListObjects.Count)
    $justAdded.Name = $name

    foreach ($lo in $ws.ListObjects) {
        if ($lo.Range.Address() -eq $ws.Range($rangeAddr).Address()) {
            $lo.TableStyle = $style
            break
        }
    }
}

New-NamedTable $ws "A14:B21" "Tabela2"  "TableStyleMedium20"
New-NamedTable $ws "C14:C22" "Tabela3"  "TableStyleMedium20"
New-NamedTable $ws "E14:F21" "Tabela26" "TableStyleMedium20"
New-NamedTable $ws "G14:G22" "Tabela6"  "TableStyleMedium20"

# --- Cosmetic view updates ------------------------------------------------

$ws.Columns("A").ColumnWidth = 9.5

[void]$ws.Range("I22").Select()
$excel.ActiveWindow.Zoom = 145

Write-Output "done"
